{"js": "// Office.js (Word JavaScript API) script.\n// Applies the benchmark-table edits described in the commit:\n//   - updates 12 single-value cells (rows 0-11, 0-indexed)\n//   - collapses the tab-separated multi-value rows 43-45 down to a single value\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// row index (0-based) -> new cell text\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"75\",\n  4: \"0.00003\",\n  5: \"0.00011\",\n  6: \"0.00006\",\n  7: \"0.00001\",\n  8: \"0.00008\",\n  9: \"0.00009\",\n  10: \"0.00009\",\n  11: \"0.00558\",\n  43: \"100\",\n  44: \"0.01\",\n  45: \"230\",\n};\n\nfor (const [row, text] of Object.entries(updates)) {\n  table.getCell(Number(row), 0).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the benchmark-table edits described in the commit:\n#   - updates 12 single-value cells (rows 1-12, 1-indexed to match Table.Cell)\n#   - collapses the tab-separated multi-value rows 44-46 down to a single value\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"75\"\n$t.Cell(5, 1).Range.Text = \"0.00003\"\n$t.Cell(6, 1).Range.Text = \"0.00011\"\n$t.Cell(7, 1).Range.Text = \"0.00006\"\n$t.Cell(8, 1).Range.Text = \"0.00001\"\n$t.Cell(9, 1).Range.Text = \"0.00008\"\n$t.Cell(10, 1).Range.Text = \"0.00009\"\n$t.Cell(11, 1).Range.Text = \"0.00009\"\n$t.Cell(12, 1).Range.Text = \"0.00558\"\n\n$t.Cell(44, 1).Range.Text = \"100\"\n$t.Cell(45, 1).Range.Text = \"0.01\"\n$t.Cell(46, 1).Range.Text = \"230\"\n"}
